# "Generate Report for Archive"
#
# 1) Status text update: every "Ready for handoff" cell (the localization
#    status shared by the Overview summary columns and the per-language
#    "Status" column) moves on to "In Translation".
# 2) The Status-related columns are narrower in the refreshed report
#    (report regeneration recomputed their auto-fit width).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1) Update status values -------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for both data rows.
$overviewStatusCells = @("E2", "F2", "E3", "F3")
foreach ($addr in $overviewStatusCells) {
    $cell = $wsOverview.Range($addr)
    if ($cell.Value2 -eq "Ready for handoff") {
        $cell.Value2 = "In Translation"
    }
}

# zh-cn / de-de sheets: "Status" column (C) for both data rows.
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    foreach ($addr in @("C2", "C3")) {
        $cell = $ws.Range($addr)
        if ($cell.Value2 -eq "Ready for handoff") {
            $cell.Value2 = "In Translation"
        }
    }
}

# --- 2) Narrow the Status columns -------------------------------------------
# Target stored column width is ~13.41 (down from ~17.22). Excel's
# ColumnWidth property only accepts values that snap to the underlying
# pixel grid, so we use the closest attainable width.
$newStatusColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth   # column F (de-de status)
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth       # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColumnWidth       # column C (Status)
